# Update "want to go" counts (column F) on both the "展览" sheet and the
# "全部类型" sheet, which duplicates the same rows.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 393
$ws1.Range("F4").Value = 3059
$ws1.Range("F6").Value = 635

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 393
$ws4.Range("F5").Value = 3059
$ws4.Range("F7").Value = 635
